$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update test-data rows: codigoTransaccion changed from 0322 to 0533 in row 2
$ws.Range("I2").Value = "0533"

# Row 3 now represents the "error" scenario for this transaction
$ws.Range("G3").Value = "Error"
$ws.Range("J3").Value = "CUENTA NO AUTORIZADA A LA SOLICITUD"
$ws.Range("H3").Value = "004"
$ws.Range("I3").Value = "0533"

# Update the active sheet view/selection to match the edited cell
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G3").Select()
